$d = $word.ActiveDocument

# Helper: replace the text of a run that starts at a given absolute
# character offset and has a known *original* text, with new text.
# Verifies the existing text matches expectations before touching it
# (so a structural mismatch fails loudly instead of corrupting text),
# and returns the length of the newly inserted text so callers can
# compute the offset of the next run.
function Set-RunText($startOffset, $expectedOldText, $newText) {
    $range = $d.Range($startOffset, $startOffset + $expectedOldText.Length)
    if ($range.Text -ne $expectedOldText) {
        throw "Run text mismatch at offset $startOffset. Expected [$expectedOldText] but found [$($range.Text)]"
    }
    $range.Text = $newText
    return $newText.Length
}

# ---------------------------------------------------------------
# Bullet 1 (paragraph 23): "Reconnaissance et sensibilisation limitées de la marque : ..."
# ---------------------------------------------------------------
$p = $d.Paragraphs.Item(23)
$pStart = $p.Range.Start

$len1 = Set-RunText $pStart "Reconnaissance et" "Reconnaissance et notoriété limitées de la marque"
$offset2 = $pStart + $len1
$oldRun2 = " sensibilisation limitées de la marque : l’obtention d’une visibilité sur ces nouveaux marchés est un obstacle majeur, nécessitant des efforts marketing robustes pour renforcer la présence de la marque Adatum à partir du terrain."
$newRun2 = " : obtenir une bonne visibilité sur ces nouveaux marchés est un défi majeur, nécessitant des initiatives marketing robustes pour établir la présence de la marque Adatum en partant de zéro."
Set-RunText $offset2 $oldRun2 $newRun2 | Out-Null

# ---------------------------------------------------------------
# Bullet 2 (paragraph 24): "Concurrence intense : ..."
# ---------------------------------------------------------------
$p = $d.Paragraphs.Item(24)
$pStart = $p.Range.Start

$len1 = Set-RunText $pStart "Concurrence" "Concurrence intense"
$offset2 = $pStart + $len1
$oldRun2 = " intense : le secteur des services cloud au Canada est très concurrentiel, avec de nombreux acteurs."
$newRun2 = " : le secteur des services cloud au Canada est très concurrentiel et comporte de nombreux acteurs."
Set-RunText $offset2 $oldRun2 $newRun2 | Out-Null

# ---------------------------------------------------------------
# Bullet 3 (paragraph 25): "Préférences et attentes diverses des clients : ..."
# ---------------------------------------------------------------
$p = $d.Paragraphs.Item(25)
$pStart = $p.Range.Start

$len1 = Set-RunText $pStart "Préférences et attentes" "Préférences et attentes diverses des clients"
$offset2 = $pStart + $len1
$oldRun2 = " diverses des clients : la mise en place de produits et de marketing adaptés aux exigences variées de ces marchés est essentielle pour la resonation avec les entreprises et les consommateurs locaux."
$newRun2 = " : il est essentiel d’adapter les produits et la campagne marketing aux exigences variées de ces marchés pour trouver un écho auprès des entreprises et des consommateurs locaux."
Set-RunText $offset2 $oldRun2 $newRun2 | Out-Null

# ---------------------------------------------------------------
# Bullet 4 (paragraph 26): "Défis en matière de réglementation et de conformité : ..."
# ---------------------------------------------------------------
$p = $d.Paragraphs.Item(26)
$pStart = $p.Range.Start

$len1 = Set-RunText $pStart "Défis" "Défis en matière de réglementation et de conformité"
$offset2 = $pStart + $len1
$oldRun2 = " réglementaires et de conformité : Adatum est confronté à la tâche complexe de naviguer dans la confidentialité, la sécurité et les réglementations opérationnelles distinctes de la région, ce qui nécessite des efforts de conformité rigoureux."
$newRun2 = " : Adatum doit identifier les règles de confidentialité et de sécurité, ainsi que les réglementations opérationnelles spécifiques de la région, une tâche complexe qui nécessite des efforts rigoureux pour respecter la conformité."
Set-RunText $offset2 $oldRun2 $newRun2 | Out-Null

# ---------------------------------------------------------------
# Bullet 5 (paragraph 27): "Complexités opérationnelles et logistiques : ..."
# Heading text is unchanged here; only the trailing description changes.
# ---------------------------------------------------------------
$p = $d.Paragraphs.Item(27)
$pStart = $p.Range.Start

$offset2 = $pStart + "Complexités opérationnelles et logistiques".Length
$oldRun2 = " : l’établissement d’opérations efficaces et interrégions présente des défis logistiques, en particulier dans le maintien de niveaux de service élevés et la gestion des centres de données dans des emplacements géographiques."
$newRun2 = " : établir des opérations efficaces et interrégions présente des défis logistiques, en particulier pour maintenir des niveaux de service élevés et gérer les centres de données dans différents emplacements géographiques."
Set-RunText $offset2 $oldRun2 $newRun2 | Out-Null
